$d = $word.ActiveDocument
$vtab = [char]11

function Replace-InRange($range, $search, $replace) {
    $f = $range.Find
    $ok = $f.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "FAILED: $search"
    }
    return $ok
}

# ---------------------------------------------------------------
# Paragraph 1: Title
# ---------------------------------------------------------------
Replace-InRange $d.Content "Quantum Computing: Unveiling Mysteries" "Unveiling the Enigma of Chemistry: The Symphony of Elements"

# ---------------------------------------------------------------
# Paragraph 2: Author line ("Dr" + "." + " Ethan Carter" -> "Amelia Stevens")
# ---------------------------------------------------------------
Replace-InRange $d.Content "Dr. Ethan Carter" "Amelia Stevens"

# ---------------------------------------------------------------
# Paragraph 3: Email line
#   "Ethan" -> "ameliasievans@emailworld" (scoped to paragraph 3 to avoid
#   touching the "Ethan" inside the author line above)
#   "Carter@QuantTech.Org" -> "com" (removes the middle runs too)
# ---------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
Replace-InRange $p3.Range "Ethan" "ameliasievans@emailworld"
$p3 = $d.Paragraphs.Item(3)
Replace-InRange $p3.Range "Carter@QuantTech.Org" "com"

# ---------------------------------------------------------------
# Paragraph 5: Body text
# ---------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)

# OP1: Replace heading run text and remove the following run (with its
# line break) that held the first "Quantum computing heralds..." sentence.
$search1 = "Quantum Computing - a Realm of Possibilities:" + $vtab + "Quantum computing heralds a new era of computation that transcends the limitations of classical computing paradigms"
$replace1 = "The world of Chemistry is a fascinating tapestry of elements and molecules, where intricate interactions orchestrate the symphony of life"
Replace-InRange $p5.Range $search1 $replace1

# OP2: Replace the "Harnessing..." sentence and extend it with two new
# sentences (merged into the same run by the engine).
$p5 = $d.Paragraphs.Item(5)
$search2 = " Harnessing the fundamental principles of quantum mechanics, this emerging field unveils the potential for groundbreaking advancements in diverse domains, ranging from optimization and cryptography to artificial intelligence and simulations"
$replace2 = " This realm of science unveils the enigmatic secrets of matter, revealing the fundamental building blocks of our universe and the processes that shape our existence. From the smallest atom to the vast expanse of galaxies, Chemistry holds the key to comprehending the delicate balance and beauty of the cosmos"
Replace-InRange $p5.Range $search2 $replace2

# OP3: Replace "In the heart of quantum computing resides the qubit..."
# (keep its leading line break untouched).
$p5 = $d.Paragraphs.Item(5)
$search3 = "In the heart of quantum computing resides the qubit, an enigmatic entity that defies the binary confinement of traditional bits"
$replace3 = "As we delve into the depths of Chemistry, we unravel the captivating artistry of chemical reactions, where elements dance in a harmonious ballet, transforming into new substances with novel properties"
Replace-InRange $p5.Range $search3 $replace3

# OP4: Replace "Qubits dance in superposition..."
$p5 = $d.Paragraphs.Item(5)
$search4 = " Qubits dance in superposition, simultaneously inhabiting states of 0, 1, or an infinite spectrum in between"
$replace4 = " The study of Chemistry empowers us to understand the intricate web of life, from the intricate workings of cellular processes to the vast array of compounds that make up the natural world"
Replace-InRange $p5.Range $search4 $replace4

# OP5: Replace "This extraordinary characteristic..."
$p5 = $d.Paragraphs.Item(5)
$search5 = " This extraordinary characteristic unleashes a computational power unfathomable by classical systems, unlocking avenues for processing that were once consigned to the realm of science fiction"
$replace5 = " It grants us the tools to unravel the mysteries of disease, paving the way for innovative treatments and therapies"
Replace-InRange $p5.Range $search5 $replace5

# OP6: Replace "While still in its nascent stages..."
# (keep its leading line break untouched).
$p5 = $d.Paragraphs.Item(5)
$search6 = "While still in its nascent stages, quantum computing has already ignited a surge of transformative applications"
$replace6 = "Furthermore, Chemistry plays a pivotal role in addressing global challenges, such as the development of sustainable energy sources, the creation of innovative materials, and the quest for cleaner and safer technologies"
Replace-InRange $p5.Range $search6 $replace6

# OP7: Replace "As we cultivate our understanding..." and remove the two
# trailing sentences that follow it.
$p5 = $d.Paragraphs.Item(5)
$search7 = " As we cultivate our understanding of this nascent technology, practical applications are emerging across industries" + "." + " From developing novel materials and optimizing supply chains to breaking cryptographic codes and accelerating drug discovery, quantum computing stands poised to revolutionize the very fabric of our world"
$replace7 = " By harnessing the power of Chemistry, we can create a sustainable future, ensuring the well-being of generations to come"
Replace-InRange $p5.Range $search7 $replace7

# ---------------------------------------------------------------
# Paragraph 7: Summary body text
# ---------------------------------------------------------------
$p7 = $d.Paragraphs.Item(7)
$searchS1 = "Quantum computing, fueled by the enigmatic power of qubits, heralds a paradigm shift in the realm of computation"
$replaceS1 = "In conclusion, Chemistry is a captivating science that unveils the mysteries of matter, unravels the intricate symphony of chemical reactions, and empowers us to address global challenges"
Replace-InRange $p7.Range $searchS1 $replaceS1

$p7 = $d.Paragraphs.Item(7)
$searchS2 = " Its implications resonate across a multitude of disciplines, spanning fields as diverse as cryptography, optimization, machine learning, and simulations" + "." + " As we delve deeper into the intricacies of this burgeoning technology, practical applications are materializing, redefining industries and transforming the world we live in"
$replaceS2 = " Its study provides a profound understanding of the universe, allowing us to harness the power of elements and molecules to create innovative solutions and shape a better world"
Replace-InRange $p7.Range $searchS2 $replaceS2

# ---------------------------------------------------------------
# Add the trailing empty paragraph that appears after the Summary
# paragraph (before the section break).
# ---------------------------------------------------------------
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange.InsertParagraphAfter()
